$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# Version 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was empty) -> Alvearie Team
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row; it must be removed,
# shifting rows 12-21 up to rows 11-20.
$wsMeta.Rows.Item(11).Delete()

# --- Elements sheet updates (row 2 = root Extension row) ---
# K2 (Short): "Extension" -> "Measure Parameter"
$wsElem.Range("K2").Value = "Measure Parameter"
# L2 (Definition): "An Extension" -> "Defined parameter options implemented by the measure"
$wsElem.Range("L2").Value = "Defined parameter options implemented by the measure"
